# Add 2022-Q3 data
# 1. Insert a new worksheet "2022-Q3" right after "总计", built from a copy of
#    the existing "2022-Q2" sheet (so it inherits identical styles/number
#    formats/column layout), then replace its single data row with the new
#    quarter's fund-holding figures.
# 2. Update the "总计" (totals) summary sheet: insert a new row for 2022-Q3
#    at the top of the data (row 2), push the existing rows down, and
#    renumber the leading index column sequentially.

$wb = $excel.ActiveWorkbook

$totals = $wb.Worksheets.Item(1)       # "总计"
$q2Sheet = $wb.Worksheets.Item(2)      # "2022-Q2" (existing, about to become the template)

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q3" sheet from a copy of "2022-Q2"
# ---------------------------------------------------------------------------
$q2Sheet.Copy($null, $totals)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Drop all the old fund rows (rows 3-9), keeping only the header (row 1)
# and a single data row (row 2) that we'll overwrite below.
$q3Sheet.Range("A3:H9").Delete()

$q3Sheet.Range("A2").Value = 0

$cell = $q3Sheet.Range("B2")
$cell.NumberFormat = "@"
$cell.Value = "001951"
$cell.ClearFormats()

$q3Sheet.Range("C2").Value = "金鹰改革红利灵活配置混合"

$cell = $q3Sheet.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.23"
$cell.ClearFormats()

$cell = $q3Sheet.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "94.49"
$cell.ClearFormats()

$cell = $q3Sheet.Range("F2")
$cell.NumberFormat = "@"
$cell.Value = "4.00"
$cell.ClearFormats()

$cell = $q3Sheet.Range("G2")
$cell.NumberFormat = "@"
$cell.Value = "1.1692"
$cell.ClearFormats()

$q3Sheet.Range("H2").Value = 8

# ---------------------------------------------------------------------------
# Step 2: update the "总计" summary sheet with the new quarter's totals
# ---------------------------------------------------------------------------
$totals.Rows.Item(2).Insert()

$totals.Range("A2:D2").ClearFormats()
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 1.17

# Restore the A-column style (bold/centered/bordered, matching the rest of
# the index column) by copying it down from the row below.
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

# Renumber the index column sequentially (0..7) now that a row was added.
for ($r = 2; $r -le 9; $r++) {
    $totals.Range("A$r").Value = $r - 2
}

# Keep the originally-active tab ("总计") selected, as in the source workbook.
$totals.Select()

